# Fix typos in the "CBM thickness" sheet:
#   - "Calson et al., 2003 (muscle)"  -> "Carlson et al., 2003 (muscle)"
#   - "Calson et al., 2003 (retina)"  -> "Carlson et al., 2003 (retina)"
#   - "Cuthbertson et al., 1986"      -> "Cuthbertson & Mandel, 1986"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM thickness")
$ws.Activate()

$ws.Range("A5").Value = "Carlson et al., 2003 (retina)"
$ws.Range("A6").Value = "Carlson et al., 2003 (muscle)"
$ws.Range("A2").Value = "Cuthbertson & Mandel, 1986"

# Match the updated cursor/selection position recorded in the saved file.
$ws.Range("A2").Select()
